$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 240.07143
$ws.Range("J2").Value = 180
$ws.Range("L2").Value = 180
$ws.Range("N2").Value = -406

# Row 40
$ws.Range("H40").Value = 1290
$ws.Range("I40").Value = 1290
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1290
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1115
$ws.Range("N40").ClearContents()

# Row 82
$ws.Range("H82").Value = 4610.154
$ws.Range("I82").Value = 2036.8889
$ws.Range("J82").Value = 10400
$ws.Range("K82").Value = 6110.6667
$ws.Range("L82").Value = 31200
$ws.Range("M82").Value = -5704.6667
$ws.Range("N82").Value = -32012

# Row 85
$ws.Range("H85").Value = 4610.154
$ws.Range("I85").Value = 2036.8889
$ws.Range("J85").Value = 10400
$ws.Range("K85").Value = 6110.6667
$ws.Range("L85").Value = 31200
$ws.Range("M85").Value = -4706.6667
$ws.Range("N85").Value = -34008

# Row 129
$ws.Range("H129").Value = 790.0476
$ws.Range("J129").Value = 943.1875
$ws.Range("L129").Value = 2829.5625
$ws.Range("N129").Value = -12829.5625

# Row 133
$ws.Range("H133").Value = 42780
$ws.Range("J133").Value = 42780
$ws.Range("L133").Value = 42780
$ws.Range("N133").Value = -52900

# Row 138
$ws.Range("H138").Value = 4263.621
$ws.Range("I138").Value = 3028.1667
$ws.Range("J138").Value = 4833.8306
$ws.Range("K138").Value = 9084.500100000001
$ws.Range("L138").Value = 14501.4918
$ws.Range("M138").Value = -3944.500100000001
$ws.Range("N138").Value = -24781.4918

# Row 141
$ws.Range("H141").Value = 4328.7334
$ws.Range("I141").Value = 4388.5386
$ws.Range("K141").Value = 13165.6158
$ws.Range("M141").Value = -7985.6158

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 89575.234
$ws.Range("I32").Value = 48381
$ws.Range("J32").Value = 102250.38
$ws.Range("K32").Value = 48381
$ws.Range("L32").Value = 102250.38
$ws.Range("M32").Value = -48094
$ws.Range("N32").Value = -102824.38

$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 28663.334
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 28663.334
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 28663.334
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -29731.334

# Row 40
$ws.Range("H40").Value = 25586
$ws.Range("J40").Value = 25586
$ws.Range("L40").Value = 25586
$ws.Range("N40").Value = -26116

# Row 75
$ws.Range("H75").Value = 10678.5
$ws.Range("I75").Value = 6485.6
$ws.Range("J75").Value = 17666.666
$ws.Range("K75").Value = 6485.6
$ws.Range("L75").Value = 17666.666
$ws.Range("M75").Value = -5549.6
$ws.Range("N75").Value = -19538.666

# Row 78
$ws.Range("H78").Value = 10678.5
$ws.Range("I78").Value = 6485.6
$ws.Range("J78").Value = 17666.666
$ws.Range("K78").Value = 19456.8
$ws.Range("L78").Value = 52999.99800000001
$ws.Range("M78").Value = -14776.8
$ws.Range("N78").Value = -62359.99800000001

# Row 93
$ws.Range("H93").Value = 31224
$ws.Range("J93").Value = 31224
$ws.Range("L93").Value = 31224
$ws.Range("N93").Value = -34968

# Row 134
$ws.Range("H134").Value = 1221.6923
$ws.Range("I134").Value = 1221.6923
$ws.Range("K134").Value = 3665.0769
$ws.Range("M134").Value = -1130.0769

$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 12803.77
$ws.Range("I41").Value = 4676.5
$ws.Range("J41").Value = 19770
$ws.Range("K41").Value = 4676.5
$ws.Range("L41").Value = 19770
$ws.Range("M41").Value = -4248.5
$ws.Range("N41").Value = -20626

# Row 52
$ws.Range("H52").Value = 55250
$ws.Range("J52").Value = 55250
$ws.Range("L52").Value = 55250
$ws.Range("N52").Value = -55838

$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Range("H60").Value = 821.875
$ws.Range("I60").Value = 113.166664
$ws.Range("J60").Value = 2948
$ws.Range("K60").Value = 339.499992
$ws.Range("L60").Value = 8844
$ws.Range("M60").Value = -88.49999200000002
$ws.Range("N60").Value = -9346

# Row 122
$ws.Range("H122").Value = 1237.9
$ws.Range("I122").Value = 513.1667
$ws.Range("J122").Value = 1721.0555
$ws.Range("K122").Value = 4618.5003
$ws.Range("L122").Value = 15489.4995
$ws.Range("M122").Value = -2168.5003
$ws.Range("N122").Value = -20389.4995

# Row 137
$ws.Range("H137").Value = 13893934
$ws.Range("I137").Value = 2300
$ws.Range("J137").Value = 16672260
$ws.Range("K137").Value = 6900
$ws.Range("L137").Value = 50016780
$ws.Range("M137").Value = -1800
$ws.Range("N137").Value = -50026980

$ws = $wb.Worksheets.Item("GSM")
# Row 137
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 45000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -55200

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3435.75
$ws.Range("I40").Value = 2926.5715
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 2926.5715
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -2790.5715
$ws.Range("N40").Value = -7272

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 27635
$ws.Range("J54").Value = 27635
$ws.Range("L54").Value = 27635
$ws.Range("N54").Value = -28675

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# Row 81
$ws.Range("H81").Value = 941.7143
$ws.Range("I81").Value = 941.7143
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1883.4286
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -822.4286
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 941.7143
$ws.Range("I84").Value = 941.7143
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9417.143
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4113.143
$ws.Range("N84").ClearContents()
